$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: insert "<comment>c_101v_02</comment>" (three differently
# formatted runs) right after the existing comment reference, before
# the run that holds "</ab>".
# ---------------------------------------------------------------------

$full = $d.Content.Text
$idx = $full.IndexOf("d'aultr</ab>")
$pos = $idx + 7   # position right after the "r" that carries the comment

# Step 1: duplicate the plain black "r" run (formatting: color 000000 +
# rtl only) right at $pos - this correctly lands AFTER the
# commentRangeEnd/commentReference markers (unlike a plain InsertBefore
# on a collapsed range, which would merge into the "r" run itself).
$srcPlain = $d.Range($pos - 1, $pos)
$srcPlain.Copy()
$tPlain = $d.Range($pos, $pos)
$tPlain.Collapse(1)
$tPlain.Paste()

# rename the pasted run's text to the comment id
$seg2 = $d.Range($pos, $pos + 1)
$seg2.Text = "c_101v_02"

# Step 2: duplicate the Courier-New / "</ab>" style run (still located
# immediately after our freshly inserted "c_101v_02" run) and paste a
# copy BEFORE "c_101v_02" -> becomes "<comment>"
$seg2EndPos = $pos + 9
$srcCourier = $d.Range($seg2EndPos, $seg2EndPos + 5)
$srcCourier.Copy()

$tBefore = $d.Range($pos, $pos)
$tBefore.Collapse(1)
$tBefore.Paste()

$seg1 = $d.Range($pos, $pos + 5)
$seg1.Text = "<comment>"
$seg1Fmt = $d.Range($pos, $pos + 9)
$seg1Fmt.Font.Color = 255          # RGB(0,0,255) -> blue

# Step 3: paste another copy of the Courier-New run AFTER "c_101v_02"
# -> becomes "</comment>"
$seg3InsertPos = $pos + 9 + 9
$tAfter = $d.Range($seg3InsertPos, $seg3InsertPos)
$tAfter.Collapse(1)
$tAfter.Paste()

$seg3 = $d.Range($seg3InsertPos, $seg3InsertPos + 5)
$seg3.Text = "</comment>"
$seg3Fmt = $d.Range($seg3InsertPos, $seg3InsertPos + 10)
$seg3Fmt.Font.Color = 255          # RGB(0,0,255) -> blue

# ---------------------------------------------------------------------
# Change 2: sectPr page margins gain a footer distance of 720 twips
# (36 pt / 0.5").
# ---------------------------------------------------------------------

$d.PageSetup.FooterDistance = 36
